$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 198.5625
$ws.Range("I5").Value = 79.5
$ws.Range("J5").Value = 397
$ws.Range("K5").Value = 79.5
$ws.Range("L5").Value = 397
$ws.Range("M5").Value = 35.5
$ws.Range("N5").Value = -627
$ws.Range("H15").Value = 1809.4584
$ws.Range("I15").Value = 1809.4584
$ws.Range("K15").Value = 5428.3752
$ws.Range("M15").Value = -5259.3752
$ws.Range("H19").Value = 47620320
$ws.Range("I19").Value = 333333340
$ws.Range("J19").Value = 1482
$ws.Range("K19").Value = 333333340
$ws.Range("L19").Value = 1482
$ws.Range("M19").Value = -333333165
$ws.Range("N19").Value = -1832
$ws.Range("H112").Value = 1431.5714
$ws.Range("I112").Value = 480
$ws.Range("J112").Value = 1479.15
$ws.Range("K112").Value = 1440
$ws.Range("L112").Value = 4437.450000000001
$ws.Range("M112").Value = -332
$ws.Range("N112").Value = -6653.450000000001
$ws.Range("H129").Value = 1039.9193
$ws.Range("J129").Value = 1103.0702
$ws.Range("L129").Value = 3309.2106
$ws.Range("N129").Value = -13309.2106
$ws.Range("H137").Value = 3074.1667
$ws.Range("I137").Value = 2411.1562
$ws.Range("J137").Value = 4400.1875
$ws.Range("K137").Value = 7233.4686
$ws.Range("L137").Value = 13200.5625
$ws.Range("M137").Value = -4683.4686
$ws.Range("N137").Value = -18300.5625
$ws.Range("H138").Value = 3193
$ws.Range("I138").Value = 1650.3462
$ws.Range("J138").Value = 4678.5186
$ws.Range("K138").Value = 4951.0386
$ws.Range("L138").Value = 14035.5558
$ws.Range("M138").Value = 188.9614000000001
$ws.Range("N138").Value = -24315.5558

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6850.84
$ws.Range("I61").Value = 3542.9143
$ws.Range("J61").Value = 14569.333
$ws.Range("K61").Value = 3542.9143
$ws.Range("L61").Value = 14569.333
$ws.Range("M61").Value = -3330.9143
$ws.Range("N61").Value = -14993.333
$ws.Range("H74").Value = 1732.4572
$ws.Range("I74").Value = 1886.0869
$ws.Range("K74").Value = 1886.0869
$ws.Range("M74").Value = -1012.0869
$ws.Range("H77").Value = 1732.4572
$ws.Range("I77").Value = 1886.0869
$ws.Range("K77").Value = 9430.434499999999
$ws.Range("M77").Value = -5062.434499999999
$ws.Range("H136").Value = 6850.84
$ws.Range("I136").Value = 3542.9143
$ws.Range("J136").Value = 14569.333
$ws.Range("K136").Value = 10628.7429
$ws.Range("L136").Value = 43707.999
$ws.Range("M136").Value = -8078.742899999999
$ws.Range("N136").Value = -48807.999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2333794.2
$ws.Range("I58").Value = 4547681
$ws.Range("K58").Value = 4547681
$ws.Range("M58").Value = -4547478
$ws.Range("H132").Value = 4209.6597
$ws.Range("I132").Value = 4266.3784
$ws.Range("J132").Value = 3999.8
$ws.Range("K132").Value = 12799.1352
$ws.Range("L132").Value = 11999.4
$ws.Range("M132").Value = -10269.1352
$ws.Range("N132").Value = -17059.4
$ws.Range("H134").Value = 2573.2456
$ws.Range("I134").Value = 1421.575
$ws.Range("J134").Value = 5283.0586
$ws.Range("K134").Value = 4264.725
$ws.Range("L134").Value = 15849.1758
$ws.Range("M134").Value = -1729.725
$ws.Range("N134").Value = -20919.1758
$ws.Range("H136").Value = 2333794.2
$ws.Range("I136").Value = 4547681
$ws.Range("K136").Value = 13643043
$ws.Range("M136").Value = -13640493

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 36.37037
$ws.Range("J2").Value = 37.384617
$ws.Range("L2").Value = 224.307702
$ws.Range("N2").Value = -450.307702
$ws.Range("H86").Value = 738
$ws.Range("I86").Value = 767.3333
$ws.Range("J86").Value = 720.4
$ws.Range("K86").Value = 2301.9999
$ws.Range("L86").Value = 2161.2
$ws.Range("M86").Value = -1115.9999
$ws.Range("N86").Value = -4533.2
$ws.Range("H89").Value = 738
$ws.Range("I89").Value = 767.3333
$ws.Range("J89").Value = 720.4
$ws.Range("K89").Value = 6905.9997
$ws.Range("L89").Value = 6483.599999999999
$ws.Range("M89").Value = -977.9997000000003
$ws.Range("N89").Value = -18339.6
$ws.Range("H97").Value = 3292.3333
$ws.Range("I97").Value = 1283.3334
$ws.Range("J97").Value = 5301.3335
$ws.Range("K97").Value = 3850.0002
$ws.Range("L97").Value = 15904.0005
$ws.Range("N97").Value = -16896.0005
$ws.Range("M97").Value = -3354.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3832.5293
$ws.Range("I102").Value = 3398.8696
$ws.Range("J102").Value = 4739.273
$ws.Range("K102").Value = 3398.8696
$ws.Range("L102").Value = 4739.273
$ws.Range("M102").Value = -1776.8696
$ws.Range("N102").Value = -7983.273
$ws.Range("H107").Value = 421.44446
$ws.Range("I107").Value = 206.46153
$ws.Range("K107").Value = 206.46153
$ws.Range("M107").Value = 1713.53847
$ws.Range("H122").Value = 18981.3
$ws.Range("I122").Value = 26076.5
$ws.Range("K122").Value = 78229.5
$ws.Range("M122").Value = -75779.5
$ws.Range("H132").Value = 4539.6924
$ws.Range("J132").Value = 19721
$ws.Range("L132").Value = 59163
$ws.Range("N132").Value = -64223

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3630.7646
$ws.Range("I7").Value = 3445.963
$ws.Range("J7").Value = 4343.5713
$ws.Range("K7").Value = 3445.963
$ws.Range("L7").Value = 4343.5713
$ws.Range("M7").Value = -3333.963
$ws.Range("N7").Value = -4567.5713
$ws.Range("H122").Value = 5784.086
$ws.Range("I122").Value = 5621.2446
$ws.Range("J122").Value = 6077.2
$ws.Range("K122").Value = 16863.7338
$ws.Range("L122").Value = 18231.6
$ws.Range("M122").Value = -14413.7338
$ws.Range("N122").Value = -23131.6
$ws.Range("H126").Value = 3630.7646
$ws.Range("I126").Value = 3445.963
$ws.Range("J126").Value = 4343.5713
$ws.Range("K126").Value = 10337.889
$ws.Range("L126").Value = 13030.7139
$ws.Range("M126").Value = -7867.889000000001
$ws.Range("N126").Value = -17970.7139
$ws.Range("H136").Value = 4276.283
$ws.Range("I136").Value = 2489.4333
$ws.Range("J136").Value = 6606.9565
$ws.Range("K136").Value = 7468.2999
$ws.Range("L136").Value = 19820.8695
$ws.Range("M136").Value = -4918.2999
$ws.Range("N136").Value = -24920.8695

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1525.421
$ws.Range("I132").Value = 747.93024
$ws.Range("J132").Value = 3913.4285
$ws.Range("K132").Value = 2243.79072
$ws.Range("L132").Value = 11740.2855
$ws.Range("M132").Value = 286.20928
$ws.Range("N132").Value = -16800.2855
$ws.Range("H136").Value = 4251
$ws.Range("I136").Value = 4248.1387
$ws.Range("K136").Value = 12744.4161
$ws.Range("M136").Value = -10194.4161
